$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values
$ws.Range("B6").Value = 12.504761904761899
$ws.Range("B7").Value = 1.2523809523809499

# New instance7 row: B8 inherits B3's old formatting (border-less, applyFont)
$ws.Range("B3").Copy()
$ws.Range("B8").PasteSpecial(-4122)

$ws.Range("A8").Value = "instance7"
$ws.Range("B8").Value = 49.884303350970001

# B3 loses its style entirely (becomes plain default)
$ws.Range("B3").ClearFormats()

# Selection moves to B11
$ws.Range("B11").Select()
